# Updates the cryptos list: refreshed prices and 1h volume percentages
# for each coin row, plus a ranking swap between Mantle and BabyDogeCoin
# (rows 45-46), as published by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.872.67'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.825.11'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6901'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07613'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.09%  '
$ws.Range("E9").Value = '  -4.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.34'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07718'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.30'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.039'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.89'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6707'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.410'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008270'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.872.18'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.078.19'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.60'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.01%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  -3.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9995'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1467'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.29'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.704'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.13'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.525'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.181'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.126'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.192'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05091'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7452'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.809'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.680'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01829'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.197.84'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.675'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9123'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.17'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9984'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.977.75'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.02%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5153'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.430'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.219'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -12.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.723'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '62.24'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -12.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4185'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.81%  '
